$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A96").Value = "2025-04-29 16:14:36"
$ws.Range("B96").Value = 283
